$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 73: drop the now-empty B73/F73/G73 cells (assigning "" clears the cell
# entirely rather than leaving an empty inline string behind).
$ws.Cells.Item(73, 2).Value = ""
$ws.Cells.Item(73, 6).Value = ""
$ws.Cells.Item(73, 7).Value = ""

# Row 74: newly scraped news item appended by the bot.
$ws.Cells.Item(74, 1).Value = '05/01/2026 11:33:04'
$ws.Cells.Item(74, 2).Value = '05/01 11:14'
$ws.Cells.Item(74, 3).Value = 'Folha de S.Paulo - Mercado - Principal'
$ws.Cells.Item(74, 4).Value = 'Banco Master: ministros do TCU acham quase impossível anular liquidação'
$ws.Cells.Item(74, 5).Value = 'https://redir.folha.com.br/redir/online/mercado/rss091/*https://www1.folha.uol.com.br/colunas/painelsa/2026/01/banco-master-ministros-do-tcu-acham-quase-impossivel-anular-liquidacao.shtml'
$ws.Cells.Item(74, 6).Value = 'tcu'
$nl = [char]10
$ws.Cells.Item(74, 7).Value = 'ia Federal, é vista por ministros do &lt;a href="https://www1.folha.uol.com.br/folha-topicos/&lt;b&gt;tcu&lt;/b&gt;/"&gt;TCU&lt;/a&gt; (Tribunal de Contas da União) como improvável.' + $nl + '&lt;a href="https://redir.folha.com'

# Assigning the multi-line value triggers Excel's implicit row auto-height
# measurement (ht/customHeight). The source row never carried an explicit
# height, so re-run AutoFit to drop the attributes again.
$ws.Rows.Item(74).AutoFit()
